# Insert a new "Data output documentation issues" slide immediately before
# the existing "Standing Questions" slide (which is currently the last
# slide in the deck).

$p = $ppt.ActivePresentation

# "Standing Questions" is the last slide right now (index == Count).
$standingQuestionsIndex = $p.Slides.Count

# Add a new slide, using the same "Title and Content" custom layout that
# the rest of the content slides use, right before it.
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$newSlide = $p.Slides.AddSlide($standingQuestionsIndex, $layout)

# --- Title -------------------------------------------------------------
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Data output documentation issues"

# --- Body / content placeholder -----------------------------------------
$body = $newSlide.Shapes.Item(2)
$tf = $body.TextFrame
$tf.AutoSize = 2

$tr = $tf.TextRange
$tr.Text = "There are a number of assays that are labeled as having been normalized by positive controls (labeled as " + [char]8220 + "blank490" + [char]8221 + "), however, the output file shows discrepancies that indicate it was not.`r[lab_id=10308,notes=NA]`r[lab_id=10309,notes=matched normal 10308]`r`r`r*these two assay outputs don" + [char]8217 + "t have a 26th column, which usually says " + [char]8220 + "blank490" + [char]8221 + " in each row. "

# Paragraphs 4 and 5 are blank separator lines with no bullet.
$para4 = $tr.Paragraphs(4,1)
$para4.ParagraphFormat.Bullet.Visible = $false
$para5 = $tr.Paragraphs(5,1)
$para5.ParagraphFormat.Bullet.Visible = $false

# Paragraph 6 is the closing note, also flush-left with no bullet.
$para6 = $tr.Paragraphs(6,1)
$para6.ParagraphFormat.Bullet.Visible = $false

# Superscript the "th" in "26th".
$para6Text = $para6.Text
$thStart = $para6.Start + $para6Text.IndexOf("26th") + 2
$thRange = $tf.TextRange.Characters($thStart, 2)
$thRange.Font.Superscript = $true
